$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 and 40: swap OKB / RenderToken content (coin order changed)
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.99"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.62"
$ws.Range("E40").Value = "  -4.48%  "

# Update Price (D) and Volume(1h) (E) values for remaining rows
$ws.Range("D2").Value = "59.405.24"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.526.23"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.99"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.62"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "2.551.40"
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("D14").Value = "2.973.78"
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.81"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "59.323.76"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "2.536.67"
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.69"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.03"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.436"
$ws.Range("E25").Value = "  -3.87%  "
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.991"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.70"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("E32").Value = "  -5.61%  "
$ws.Range("E33").Value = "  +5.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "159.99"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.75"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.42"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  -5.95%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "294.24"
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  +1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.45"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0228"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("E51").Value = "  -3.02%  "
